# Swap data between row 3 and row 4 for columns A, B, D, E, F, G, H, I
# (columns C and J..AY are identical between the two rows and remain unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "I")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")

    $val3 = $cell3.Value2
    $val4 = $cell4.Value2

    $cell3.Value2 = $val4
    $cell4.Value2 = $val3
}
